$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 235
    3  = 237
    4  = 240
    5  = 241
    6  = 243
    7  = 244
    8  = 246
    9  = 248
    10 = 250
    11 = 251
    12 = 253
    13 = 256
    14 = 257
    15 = 258
    16 = 0
    17 = 39
    18 = 60
    19 = 122
    20 = 134
    21 = 275
    22 = 308
    23 = 323
    24 = 343
    25 = 441
    26 = 470
    27 = 513
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
